$wb = $excel.ActiveWorkbook

$wsReview = $wb.Worksheets.Item("LH_Review_WF_LOGIN")
$wsHistory = $wb.Worksheets.Item("VERSION-HISTORY")

# Update "Owner Status" column (H) for rows 2-4 from "open" to "not applicable"
$wsReview.Range("H2").Value = "not applicable"
$wsReview.Range("H3").Value = "not applicable"
$wsReview.Range("H4").Value = "not applicable"

# Update selections to match final saved state
$wsReview.Range("I11").Select()
$wsHistory.Range("C19").Select()

# Make the review sheet the active one (tabSelected)
$wsReview.Activate()
$wsReview.Range("I11").Select()
